# Apply stock-report corrections: reduce quantities (col F) and recompute
# values (col G = Rate*Qty) for affected line items, update the company-level
# Sub Total rows (col B), the two swapped item-code rows (258/259 and 364/365),
# and the grand totals at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(98, 6).Value = 30
$ws.Cells.Item(98, 7).Value = 5161.2
$ws.Cells.Item(133, 2).Value = 205812.19
$ws.Cells.Item(167, 6).Value = 48
$ws.Cells.Item(167, 7).Value = 2031.36
$ws.Cells.Item(176, 2).Value = 14579.25
$ws.Cells.Item(211, 6).Value = 28
$ws.Cells.Item(211, 7).Value = 1181.88
$ws.Cells.Item(216, 2).Value = 9691.110000000001
$ws.Cells.Item(225, 6).Value = 16
$ws.Cells.Item(225, 7).Value = 1300.48
$ws.Cells.Item(228, 2).Value = 11503.29
$ws.Cells.Item(296, 6).Value = 3
$ws.Cells.Item(296, 7).Value = 816.36
$ws.Cells.Item(310, 2).Value = 107928.53
$ws.Cells.Item(314, 2).Value = 57077
$ws.Cells.Item(314, 4).Value = 93.08
$ws.Cells.Item(314, 5).Value = 111.2
$ws.Cells.Item(314, 6).Value = 1
$ws.Cells.Item(314, 7).Value = 93.08
$ws.Cells.Item(315, 2).Value = 61610
$ws.Cells.Item(315, 4).Value = 102.71
$ws.Cells.Item(315, 5).Value = 122.71
$ws.Cells.Item(315, 6).Value = 88
$ws.Cells.Item(315, 7).Value = 9038.48
$ws.Cells.Item(342, 6).Value = 133
$ws.Cells.Item(342, 7).Value = 6933.29
$ws.Cells.Item(351, 6).Value = 2
$ws.Cells.Item(351, 7).Value = 177.6
$ws.Cells.Item(356, 6).Value = 10
$ws.Cells.Item(356, 7).Value = 1985.9
$ws.Cells.Item(370, 6).Value = 15
$ws.Cells.Item(370, 7).Value = 2996.55
$ws.Cells.Item(380, 2).Value = 250741.58
$ws.Cells.Item(422, 6).Value = 1
$ws.Cells.Item(422, 7).Value = 350.24
$ws.Cells.Item(429, 2).Value = 6636.48
$ws.Cells.Item(436, 2).Value = 47097
$ws.Cells.Item(436, 4).Value = 112.28
$ws.Cells.Item(436, 5).Value = 134.16
$ws.Cells.Item(436, 6).Value = 15
$ws.Cells.Item(436, 7).Value = 1684.2
$ws.Cells.Item(437, 2).Value = 58047
$ws.Cells.Item(437, 4).Value = 105.54
$ws.Cells.Item(437, 5).Value = 126.1
$ws.Cells.Item(437, 6).Value = 62
$ws.Cells.Item(437, 7).Value = 6543.48
$ws.Cells.Item(497, 6).Value = 278
$ws.Cells.Item(497, 7).Value = 3561.18
$ws.Cells.Item(501, 6).Value = 83
$ws.Cells.Item(501, 7).Value = 1615.18
$ws.Cells.Item(509, 2).Value = 93704.99000000001
$ws.Cells.Item(512, 6).Value = 27
$ws.Cells.Item(512, 7).Value = 783.54
$ws.Cells.Item(516, 2).Value = 6535.16
$ws.Cells.Item(559, 6).Value = 306
$ws.Cells.Item(559, 7).Value = 2050.2
$ws.Cells.Item(560, 6).Value = 294
$ws.Cells.Item(560, 7).Value = 4859.82
$ws.Cells.Item(563, 2).Value = 36471.78
$ws.Cells.Item(608, 6).Value = 11
$ws.Cells.Item(608, 7).Value = 439.34
$ws.Cells.Item(613, 2).Value = 5940.7
$ws.Cells.Item(636, 6).Value = 44
$ws.Cells.Item(636, 7).Value = 5404.52
$ws.Cells.Item(640, 2).Value = 207609.45
$ws.Cells.Item(646, 6).Value = 5
$ws.Cells.Item(646, 7).Value = 136
$ws.Cells.Item(649, 2).Value = 52992.13
$ws.Cells.Item(652, 6).Value = 2
$ws.Cells.Item(652, 7).Value = 7231.8
$ws.Cells.Item(666, 2).Value = 36368.83
$ws.Cells.Item(668, 6).Value = 6
$ws.Cells.Item(668, 7).Value = 198.66
$ws.Cells.Item(670, 6).Value = 70
$ws.Cells.Item(670, 7).Value = 2317.7
$ws.Cells.Item(674, 6).Value = 14
$ws.Cells.Item(674, 7).Value = 463.54
$ws.Cells.Item(677, 2).Value = 20403.24
$ws.Cells.Item(685, 6).Value = 26
$ws.Cells.Item(685, 7).Value = 1627.86
$ws.Cells.Item(687, 6).Value = 26
$ws.Cells.Item(687, 7).Value = 3031.08
$ws.Cells.Item(690, 6).Value = 32
$ws.Cells.Item(690, 7).Value = 3385.28
$ws.Cells.Item(692, 6).Value = 12
$ws.Cells.Item(692, 7).Value = 1062.24
$ws.Cells.Item(694, 2).Value = 35106.83
$ws.Cells.Item(749, 6).Value = 13
$ws.Cells.Item(749, 7).Value = 7435.61
$ws.Cells.Item(752, 2).Value = 14116.2
$ws.Cells.Item(763, 6).Value = 365
$ws.Cells.Item(763, 7).Value = 13607.2
$ws.Cells.Item(770, 6).Value = 475
$ws.Cells.Item(770, 7).Value = 64129.75
$ws.Cells.Item(772, 6).Value = 543
$ws.Cells.Item(772, 7).Value = 65545.53
$ws.Cells.Item(774, 2).Value = 240636.01
$ws.Cells.Item(783, 6).Value = 179
$ws.Cells.Item(783, 7).Value = 5391.48
$ws.Cells.Item(792, 2).Value = 13312.97
$ws.Cells.Item(848, 6).Value = 2
$ws.Cells.Item(848, 7).Value = 59.84
$ws.Cells.Item(849, 2).Value = 1452.68
$ws.Cells.Item(851, 6).Value = 578
$ws.Cells.Item(851, 7).Value = 17472.94
$ws.Cells.Item(852, 6).Value = 3099
$ws.Cells.Item(852, 7).Value = 505477.89
$ws.Cells.Item(854, 6).Value = 222
$ws.Cells.Item(854, 7).Value = 32112.3
$ws.Cells.Item(855, 6).Value = 116
$ws.Cells.Item(855, 7).Value = 4424.24
$ws.Cells.Item(860, 2).Value = 617924.03
$ws.Cells.Item(866, 2).Value = 3431326.9
$ws.Cells.Item(867, 2).Value = 3431326.9
